$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row (row 12) with the new test mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Wat zijn de verzendkosten?"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("C12").Value = "Testmail #3: Wat zijn de verzendkosten?"
$logs.Range("D12").Value = "Productinformatie"
$logs.Range("F12").Value = "2025-06-29 14:15:40"
$logs.Range("G12").Value = "Nee"
$logs.Range("H12").Value = "Ja"
$logs.Range("I12").Value = "Nee"

# --- Extend the conditional formatting ranges so they cover the new row ---
$dCfs = $logs.Range("D2:D11").FormatConditions
for ($i = 1; $i -le $dCfs.Count; $i++) {
    $dCfs.Item($i).ModifyAppliesToRange($logs.Range("D2:D12"))
}

$gCfs = $logs.Range("G2:G11").FormatConditions
for ($i = 1; $i -le $gCfs.Count; $i++) {
    $gCfs.Item($i).ModifyAppliesToRange($logs.Range("G2:G12"))
}

$hCfs = $logs.Range("H2:H11").FormatConditions
for ($i = 1; $i -le $hCfs.Count; $i++) {
    $hCfs.Item($i).ModifyAppliesToRange($logs.Range("H2:H12"))
}

$iCfs = $logs.Range("I2:I11").FormatConditions
for ($i = 1; $i -le $iCfs.Count; $i++) {
    $iCfs.Item($i).ModifyAppliesToRange($logs.Range("I2:I12"))
}

# --- Sheet "Dashboard": update the Productinformatie count from 2 to 3 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B4").Value = 3
